# Auto-generated edit script: refresh cryptos price/volume snapshot
# (matches commit "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.889.17"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "2.225.53"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.48"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.24"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +9.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0969"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.36"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +8.10%  "
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "2.557.31"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.04"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "2.212.30"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "41.889.44"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +12.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +18.65%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.53"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.57"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0742"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.70"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.06"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0308"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.39%  "
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.48"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +23.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.76"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.208"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.77"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.68"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.63%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.07%  "
